$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 154.4054388999939

$ws.Range("F2").Value = "durations on validate set"

$ws.Range("A4").Value = 18814.897
$ws.Range("B4").Value = 18339

$ws.Range("F4").Value = 5945.339
$ws.Range("G4").Value = 5873
$ws.Range("H4").Value = 5856
$ws.Range("I4").Value = 5986
$ws.Range("J4").Value = 5895
